$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Merge title runs "PROJECT AND " + "PORTFOLIO PEER REVIEWS"
Replace-Text "PROJECT AND PORTFOLIO PEER REVIEWS" "PROJECT AND PORTFOLIO PEER REVIEWS"

# 2. Merge "P" + "ortfolio " runs
Replace-Text "Portfolio (HW2)" "Portfolio (HW2)"

# 3. Merge "submissions " + "are available as " + "Piazza posts. Choose " runs
Replace-Text "submissions are available as Piazza posts. Choose " "submissions are available as Piazza posts. Choose "

# 4. Remove " (.pdf)" from "Separately, upload a document (.pdf) with the text "
Replace-Text "Separately, upload a document (.pdf) with the text " "Separately, upload a document with the text "

# 5. Merge the peer-review sentence runs and drop the comma in "anonymously, but"
Replace-Text "You may post your Piazza review anonymously, but be aware" "You may post your Piazza review anonymously but be aware"

Write-Output "done"
